$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2110
$ws1.Range("F4").Value = 869
$ws1.Range("F5").Value = 1310
$ws1.Range("F6").Value = 365

# Sheet "全部类型": update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2110
$ws4.Range("F6").Value = 869
$ws4.Range("F7").Value = 1310
$ws4.Range("F8").Value = 365
